# Apply the commit's changes (getWatershed switched to the IEc file; the
# downstream tablesOut / tablesToExcel exports were regenerated):
#  - "Lower Mortendad" watershed values collapse to "Mortendad"
#  - Several Max/Last Cr values, their sample dates, and the computed
#    "Length" column were recalculated
#  - The exhibit's "Mortendad Canyon" section heading becomes
#    "Los Alamos and Pajarito Canyons"
#  - The R-19 S2 well row is no longer emitted and is removed from both
#    sheets
#
# NOTE: values such as "10.0", "50", "2008-02-20" are stored in the source
# workbook as plain TEXT (t="inlineStr"), not numbers/dates. Assigning a
# bare numeric- or date-looking string via .Value triggers Excel's normal
# "smart" type inference and would silently turn them into a number/date
# cell. Prefixing with a leading apostrophe forces a literal text entry,
# matching the original file's cell typing.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Intermediate for Mapping")
$ws2 = $wb.Worksheets.Item("Intermediate Exhibit")

# ----- Sheet 1: "Intermediate for Mapping" -----

# Row 4 (R-12 S1)
$ws1.Range("V4").Value = "'10.0"
$ws1.Range("W4").Value = "'2008-02-20"
$ws1.Range("AD4").Value = "'50"

# Row 5 (MCOBT-4.4)
$ws1.Range("P5").Value = "Mortendad"

# Row 6 (MCOI-4)
$ws1.Range("P6").Value = "Mortendad"

# Row 7 (MCOI-8)
$ws1.Range("P7").Value = "Mortendad"

# Row 8 (R-55i)
$ws1.Range("P8").Value = "Mortendad"
$ws1.Range("V8").Value = "'10.0"
$ws1.Range("W8").Value = "'2011-05-10"
$ws1.Range("X8").Value = "'10.0"
$ws1.Range("AD8").Value = "'31"

# Row 12 (R-9i S2)
$ws1.Range("X12").Value = "'10.0"
$ws1.Range("Y12").Value = "'2013-08-08"
$ws1.Range("AD12").Value = "'29"

# Row 14 (Test Well 2A)
$ws1.Range("V14").Value = "'10.0"
$ws1.Range("X14").Value = "'10.0"
$ws1.Range("AD14").Value = "'8"

# Row 15 (03-B-10)
$ws1.Range("AD15").Value = "'28"

# Row 16 (R-19 S2) no longer present - remove the entire row
$ws1.Rows.Item(16).Delete()

# ----- Sheet 2: "Intermediate Exhibit" -----

# Row 6 (R-12 S1)
$ws2.Range("G6").Value = "'10"
$ws2.Range("H6").Value = "'2/20/08"

# Row 11 (R-55i)
$ws2.Range("G11").Value = "'10"
$ws2.Range("H11").Value = "'5/10/11"

# Row 12 section heading (was "Mortendad Canyon")
$ws2.Range("A12").Value = "Los Alamos and Pajarito Canyons"

# Row 18 (Test Well 2A)
$ws2.Range("G18").Value = "'10"

# Row 20 (R-19 S2) no longer present - remove the entire row
$ws2.Rows.Item(20).Delete()
